$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same style as
# the existing header row (bold, bordered, centered) -- copy the format
# from H1 (an existing header cell) before writing the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-36.
$data = @(
    @(5, 5),
    @(7, 8),
    @(5, 6),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(9, 9),
    @(7, 7),
    @(10, 10),
    @(7, 7),
    @(8, 9),
    @(6, 7),
    @(10, 10),
    @(5, 5),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(7, 7),
    @(4, 5),
    @(6, 8),
    @(6, 7),
    @(4, 6),
    @(1, 3),
    @(1, 3),
    @(1, 3),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
